$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters -> numbers: D=4, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20

# Row 3
$ws.Cells.Item(3, 4).Value  = 44511
$ws.Cells.Item(3, 13).Value = 15
$ws.Cells.Item(3, 14).Value = 22000
$ws.Cells.Item(3, 15).Value = 22000
$ws.Cells.Item(3, 16).Value = 22000
$ws.Cells.Item(3, 19).Value = 1467

# Row 4
$ws.Cells.Item(4, 4).Value  = 44264
$ws.Cells.Item(4, 12).Value = "Calibre 100"
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 20000
$ws.Cells.Item(4, 15).Value = 20000
$ws.Cells.Item(4, 16).Value = 20000
$ws.Cells.Item(4, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(4, 19).Value = 1111
$ws.Cells.Item(4, 20).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value  = 44217
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 55
$ws.Cells.Item(5, 14).Value = 18000
$ws.Cells.Item(5, 15).Value = 18000
$ws.Cells.Item(5, 16).Value = 18000
$ws.Cells.Item(5, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(5, 19).Value = 1000
$ws.Cells.Item(5, 20).Value = 18

# Row 6
$ws.Cells.Item(6, 4).Value  = 44601
$ws.Cells.Item(6, 13).Value = 30
$ws.Cells.Item(6, 14).Value = 28000
$ws.Cells.Item(6, 15).Value = 28000
$ws.Cells.Item(6, 16).Value = 28000
$ws.Cells.Item(6, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(6, 19).Value = 1556
$ws.Cells.Item(6, 20).Value = 18

# Row 7
$ws.Cells.Item(7, 4).Value  = 44495
$ws.Cells.Item(7, 13).Value = 50
$ws.Cells.Item(7, 14).Value = 24000
$ws.Cells.Item(7, 15).Value = 24000
$ws.Cells.Item(7, 16).Value = 24000
$ws.Cells.Item(7, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(7, 18).Value = "China"
$ws.Cells.Item(7, 19).Value = 2400
$ws.Cells.Item(7, 20).Value = 10

# Row 8
$ws.Cells.Item(8, 4).Value  = 44427
$ws.Cells.Item(8, 13).Value = 55
$ws.Cells.Item(8, 14).Value = 7000
$ws.Cells.Item(8, 15).Value = 7000
$ws.Cells.Item(8, 16).Value = 7000
$ws.Cells.Item(8, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 467
$ws.Cells.Item(8, 20).Value = 15

# Row 9
$ws.Cells.Item(9, 4).Value  = 44208
$ws.Cells.Item(9, 12).Value = "Especial"
$ws.Cells.Item(9, 13).Value = 70
$ws.Cells.Item(9, 14).Value = 24000
$ws.Cells.Item(9, 15).Value = 24000
$ws.Cells.Item(9, 16).Value = 24000
$ws.Cells.Item(9, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(9, 19).Value = 1600
$ws.Cells.Item(9, 20).Value = 15

# Row 10
$ws.Cells.Item(10, 4).Value  = 44411
$ws.Cells.Item(10, 13).Value = 210
$ws.Cells.Item(10, 14).Value = 8000
$ws.Cells.Item(10, 15).Value = 8000
$ws.Cells.Item(10, 16).Value = 8000
$ws.Cells.Item(10, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(10, 19).Value = 1000
$ws.Cells.Item(10, 20).Value = 8
